# Updated main GSC export data:
# The exported date range rolled forward by one day - the oldest day
# (2025-10-26, the first data row) dropped off the front of the report
# and all remaining rows shift up by one, so the trailing row
# (previously 2026-01-20) is no longer present.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the obsolete first data row (row 2, "2025-10-26"). Excel's
# native row delete shifts every row below it up by one and shrinks the
# sheet's used range accordingly - exactly matching a rolling GSC export.
$ws.Rows.Item(2).Delete()
